# Apply the updated crypto price/volume(1h) values for this data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cellRef, $text) {
    # Leading apostrophe forces Excel to store the value as text (avoids
    # numeric/date auto-conversion of values like '510.93' or '0.998');
    # resetting the Style back to Normal afterwards strips the resulting
    # quote-prefix formatting so the cell keeps its original (default) style.
    $cell = $ws.Range($cellRef)
    $cell.Value = "'" + $text
    $cell.Style = "Normal"
}

Set-TextCell "D2" '60.366.08'
Set-TextCell "E2" '  -0.30%  '
Set-TextCell "D3" '2.608.46'
Set-TextCell "E3" '  -1.46%  '
Set-TextCell "E4" '  -0.17%  '
Set-TextCell "D5" '510.93'
Set-TextCell "E5" '  -0.30%  '
Set-TextCell "D6" '154.85'
Set-TextCell "E6" '  -1.76%  '
Set-TextCell "D7" '0.998'
Set-TextCell "E7" '  +0.47%  '
Set-TextCell "D8" '0.589'
Set-TextCell "E8" '  -2.16%  '
Set-TextCell "D9" '2.619.16'
Set-TextCell "E9" '  -2.82%  '
Set-TextCell "D10" '6.68'
Set-TextCell "E10" '  +2.74%  '
Set-TextCell "E11" '  -0.67%  '
Set-TextCell "E12" '  -0.65%  '
Set-TextCell "E13" '  +1.56%  '
Set-TextCell "D14" '3.064.81'
Set-TextCell "E14" '  -2.35%  '
Set-TextCell "D15" '60.356.30'
Set-TextCell "E15" '  -0.49%  '
Set-TextCell "D16" '21.62'
Set-TextCell "E16" '  -1.04%  '
Set-TextCell "E17" '  +0.04%  '
Set-TextCell "D18" '2.616.61'
Set-TextCell "E18" '  -2.48%  '
Set-TextCell "E19" '  -0.90%  '
Set-TextCell "D20" '350.91'
Set-TextCell "E20" '  +0.39%  '
Set-TextCell "E21" '  +0.37%  '
Set-TextCell "D22" '6.15'
Set-TextCell "D23" '0.998'
Set-TextCell "E23" '  +0.12%  '
Set-TextCell "D24" '60.50'
Set-TextCell "E24" '  +0.23%  '
Set-TextCell "D25" '0.422'
Set-TextCell "E25" '  -0.44%  '
Set-TextCell "E26" '  -0.93%  '
Set-TextCell "E27" '  +0.16%  '
Set-TextCell "D28" '0.0₃0841'
Set-TextCell "E28" '  -3.74%  '
Set-TextCell "E29" '  -2.65%  '
Set-TextCell "E30" '  +0.32%  '
Set-TextCell "E31" '  -1.30%  '
Set-TextCell "D32" '151.08'
Set-TextCell "E32" '  -4.01%  '
Set-TextCell "D33" '1.57'
Set-TextCell "E33" '  -0.76%  '
Set-TextCell "E34" '  +0.12%  '
Set-TextCell "E35" '  -1.94%  '
Set-TextCell "E36" '  -2.82%  '
Set-TextCell "D37" '0.883'
Set-TextCell "E37" '  +5.48%  '
Set-TextCell "E38" '  -1.98%  '
Set-TextCell "D39" '0.845'
Set-TextCell "E39" '  -1.38%  '
Set-TextCell "D40" '36.29'
Set-TextCell "E40" '  +2.25%  '
Set-TextCell "E41" '  -0.49%  '
Set-TextCell "D42" '295.06'
Set-TextCell "E42" '  -6.11%  '
Set-TextCell "E43" '  -3.80%  '
Set-TextCell "D44" '0.101'
Set-TextCell "E44" '  -0.11%  '
Set-TextCell "E45" '  +0.53%  '
Set-TextCell "D46" '0.0556'
Set-TextCell "E46" '  -3.68%  '
Set-TextCell "D47" '19.88'
Set-TextCell "E47" '  -1.78%  '
Set-TextCell "D48" '4.80'
Set-TextCell "E48" '  -2.14%  '
Set-TextCell "E49" '  -1.19%  '
Set-TextCell "E50" '  +0.15%  '
Set-TextCell "D51" '2.002.76'
Set-TextCell "E51" '  -3.57%  '
